# Append 5 new arrival rows (rows 70-74) to the "Main Data" sheet,
# matching the data scraped for the Jan 08 (Sunday) late-night arrivals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=70; A=69; B="Sunday, Jan 08"; C="9:20 PM";  D="FR6390"; E="Dortmund"; F="(DTM)"; G="Ryanair ";       H="B738"; I="(SP-RSO)"; J="9:25 PM";          L="0 hours, 5 minutes" },
    @{ Row=71; A=70; B="Sunday, Jan 08"; C="9:40 PM";  D="FR6869"; E="Cologne";  F="(CGN)"; G="Ryanair ";       H="B738"; I="(9H-QBA)"; J="9:38 PM";          L="0 hours, -2 minutes" },
    @{ Row=72; A=71; B="Sunday, Jan 08"; C="10:00 PM"; D="W61094"; E="Dortmund"; F="(DTM)"; G="Wizz Air ";      H="A321"; I="(HA-LXP)"; J="9:36 PM";          L="0 hours, -24 minutes" },
    @{ Row=73; A=72; B="Sunday, Jan 08"; C="10:31 PM"; D="3V4639"; E="Paris";    F="(CDG)"; G="ASL Airlines ";  H="B738"; I="(OE-IXB)"; J="Diverted to BRQ";  L=$null },
    @{ Row=74; A=73; B="Sunday, Jan 08"; C="10:33 PM"; D="3V4463"; E="Paris";    F="(CDG)"; G="FedEx ";         H="B738"; I="(OE-IWF)"; J="10:31 PM";         L="0 hours, -2 minutes" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    if ($r.L) {
        $ws.Range("L$row").Value = $r.L
    }
}
